$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update changed cell values (rows 25-52 of component errors sheet)
$ws.Range("J25").Value = -17.38961135916213
$ws.Range("K25").Value = -3.79150938822805
$ws.Range("I26").Value = -17.34253003720006
$ws.Range("J26").Value = -3.744557841334466
$ws.Range("H27").Value = -17.38776996662565
$ws.Range("I27").Value = -3.78966448320287
$ws.Range("G28").Value = -17.34253003720006
$ws.Range("H28").Value = -3.744557841334466
$ws.Range("F29").Value = -17.34740415032606
$ws.Range("G29").Value = -3.749049852560759
$ws.Range("H29").Value = -1.092443334187166
$ws.Range("I29").Value = 0.2166319823593028
$ws.Range("J29").Value = 1.738460323404823
$ws.Range("K29").Value = -4.026811231704187
$ws.Range("E30").Value = -17.34820405356395
$ws.Range("F30").Value = -3.749852217274523
$ws.Range("G30").Value = -1.091540348589674
$ws.Range("H30").Value = 0.2165081321257674
$ws.Range("I30").Value = 1.738336059163075
$ws.Range("J30").Value = -4.026935910179134
$ws.Range("D31").Value = -15.16609590571279
$ws.Range("E31").Value = -2.929546394552062
$ws.Range("F31").Value = -0.9697609669912592
$ws.Range("G31").Value = 0.2176155739527843
$ws.Range("H31").Value = 1.739447192114724
$ws.Range("I31").Value = -4.025821084199961
$ws.Range("C32").Value = -5.042530037200061
$ws.Range("D32").Value = 5.855442158665534
$ws.Range("E32").Value = 2.750960042993063
$ws.Range("F32").Value = 1.26020908108197
$ws.Range("G32").Value = 1.78156889852946
$ws.Range("H32").Value = -3.984172743600905
$ws.Range("B33").Value = -5.717024052945431
$ws.Range("C33").Value = 1.09324138922065
$ws.Range("D33").Value = 1.632021407165425
$ws.Range("E33").Value = 1.641771583357695
$ws.Range("F33").Value = 2.664585651028034
$ws.Range("G33").Value = -3.562788079019981
$ws.Range("H33").Value = 2.992524635168522
$ws.Range("I33").Value = 0.5244980819872467
$ws.Range("J33").Value = -1.226968798656231
$ws.Range("K33").Value = 1.712812639041644
$ws.Range("B34").Value = -1.731369448419173
$ws.Range("C34").Value = -1.784697304814003
$ws.Range("D34").Value = 2.133859601623712
$ws.Range("E34").Value = 3.944422896034453
$ws.Range("F34").Value = -2.888251961988018
$ws.Range("G34").Value = 3.369652086375609
$ws.Range("H34").Value = 0.7030450978900973
$ws.Range("I34").Value = -1.170280945219361
$ws.Range("J34").Value = 1.769389880649555
$ws.Range("B35").Value = -0.04213095360905994
$ws.Range("C35").Value = 1.048810608413148
$ws.Range("D35").Value = 3.434801070000262
$ws.Range("E35").Value = -2.985411091539163
$ws.Range("F35").Value = 3.363320785887396
$ws.Range("G35").Value = 0.7011015178531892
$ws.Range("H35").Value = -1.115345921577628
$ws.Range("I35").Value = 1.768116730534516
$ws.Range("B36").Value = 0.5474942774927034
$ws.Range("C36").Value = 3.594648531973121
$ws.Range("D36").Value = -2.897671109164335
$ws.Range("E36").Value = 3.386630971890952
$ws.Range("F36").Value = 0.7030914489441358
$ws.Range("G36").Value = -1.112872643084927
$ws.Range("H36").Value = 1.769646039288162
$ws.Range("B37").Value = 2.432348187196169
$ws.Range("C37").Value = -2.737865361790014
$ws.Range("D37").Value = 4.185741306446999
$ws.Range("E37").Value = 1.293909740146092
$ws.Range("F37").Value = -1.08737519398072
$ws.Range("G37").Value = 1.804298595283711
$ws.Range("H37").Value = 0.4244738359701329
$ws.Range("I37").Value = 1.930305569530376
$ws.Range("J37").Value = 1.603972340061787
$ws.Range("K37").Value = 2.442460839738445
$ws.Range("B38").Value = -4.415541225076481
$ws.Range("C38").Value = 2.949137148088397
$ws.Range("D38").Value = 3.310659458692954
$ws.Range("E38").Value = 0.7766524156165242
$ws.Range("F38").Value = 2.957275387081722
$ws.Range("G38").Value = 0.8684929633549106
$ws.Range("H38").Value = 1.928862871635673
$ws.Range("I38").Value = 1.602525255517833
$ws.Range("J38").Value = 2.441009366320828
$ws.Range("B39").Value = 0.1018375773526259
$ws.Range("C39").Value = 1.232393410787807
$ws.Range("D39").Value = 0.1773460743341531
$ws.Range("E39").Value = 2.208017199995609
$ws.Range("F39").Value = 0.5935582784456845
$ws.Range("G39").Value = 1.934518833886514
$ws.Range("H39").Value = 1.608198178209705
$ws.Range("I39").Value = 2.446699256038415
$ws.Range("B40").Value = 0.8313796242911025
$ws.Range("C40").Value = 0.01906042597315105
$ws.Range("D40").Value = 2.957275387081737
$ws.Range("E40").Value = 1.000974305187484
$ws.Range("F40").Value = 2.133050092891608
$ws.Range("G40").Value = 1.606482930593387
$ws.Range("H40").Value = 2.444978945836923
$ws.Range("B41").Value = -1.853083536387885
$ws.Range("C41").Value = 2.060743096760248
$ws.Range("D41").Value = -0.1474690440171003
$ws.Range("E41").Value = 1.473124225205538
$ws.Range("F41").Value = 2.526843383637697
$ws.Range("G41").Value = 2.865510085538972
$ws.Range("H41").Value = 0.01173936477690107
$ws.Range("I41").Value = 0.7971074916902978
$ws.Range("J41").Value = 2.521799915085936
$ws.Range("K41").Value = 2.888634420078432
$ws.Range("B42").Value = 0.5312983411637617
$ws.Range("C42").Value = 0.1836660922966331
$ws.Range("D42").Value = 2.173098942530558
$ws.Range("E42").Value = 1.92640009565342
$ws.Range("F42").Value = 2.604449907677293
$ws.Range("G42").Value = -0.3140189780649787
$ws.Range("H42").Value = 0.8018524083122998
$ws.Range("I42").Value = 2.526559023311947
$ws.Range("J42").Value = 2.893407724866548
$ws.Range("B43").Value = -1.059613776819873
$ws.Range("C43").Value = 2.064435256157637
$ws.Range("D43").Value = 2.167802532234504
$ws.Range("E43").Value = 2.84544118798722
$ws.Range("F43").Value = -0.07669164863153788
$ws.Range("G43").Value = 0.9997159764915864
$ws.Range("H43").Value = 2.52373924162616
$ws.Range("I43").Value = 2.890579586181772
$ws.Range("B44").Value = 1.033156903866584
$ws.Range("C44").Value = 1.479261134386391
$ws.Range("D44").Value = 2.823892576633995
$ws.Range("E44").Value = 0.4850868753334661
$ws.Range("F44").Value = 1.460902994935779
$ws.Range("G44").Value = 2.944112461481595
$ws.Range("H44").Value = 3.14983302187899
$ws.Range("B45").Value = 1.122595510644175
$ws.Range("C45").Value = 2.173440877464231
$ws.Range("D45").Value = -0.3391222696877914
$ws.Range("E45").Value = 1.15019388917726
$ws.Range("F45").Value = 2.825176478245041
$ws.Range("G45").Value = 3.151877404826194
$ws.Range("H45").Value = -2.471740390448148
$ws.Range("I45").Value = 0.8884784691348386
$ws.Range("B46").Value = 0.570551512911309
$ws.Range("C46").Value = 0.07843845383020209
$ws.Range("D46").Value = 1.388273133551422
$ws.Range("E46").Value = 2.905385746929681
$ws.Range("F46").Value = 3.232175349920908
$ws.Range("G46").Value = -2.392757661014926
$ws.Range("H46").Value = 0.968054944197732
$ws.Range("B47").Value = -1.340014270969163
$ws.Range("C47").Value = 0.6280135122058965
$ws.Range("D47").Value = 2.905983602396475
$ws.Range("E47").Value = 3.232973181635046
$ws.Range("F47").Value = -2.389481588625884
$ws.Range("G47").Value = 0.9697495043355957
$ws.Range("B48").Value = 1.716144675562305
$ws.Range("C48").Value = 3.240715019062477
$ws.Range("D48").Value = 3.629692481104868
$ws.Range("E48").Value = -2.485016276448988
$ws.Range("F48").Value = 0.409012827602723
$ws.Range("B49").Value = 1.626942668542526
$ws.Range("C49").Value = 2.724657014426995
$ws.Range("D49").Value = -2.773702763948436
$ws.Range("E49").Value = 0.6347777823444005
$ws.Range("B50").Value = 2.627648411919878
$ws.Range("C50").Value = -2.893259719291934
$ws.Range("D50").Value = 0.51982466486254
$ws.Range("B51").Value = -3.929515313798774
$ws.Range("C51").Value = 0.08810208718128365
$ws.Range("B52").Value = -0.1873605988422895

# Clear cells that are no longer part of the data range
$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
